# Atualizacao de bases das ligas, do dia: 17-06-2024 as 21:10
# The two matches played on the same date/time had their data rows
# swapped (id + HomeTeam..PL_AhUnder -> columns B and E:AD), while the
# row-index column (A), Div (C) and Date (D) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2) {
    $b1 = $ws.Range("B${row1}").Value2
    $b2 = $ws.Range("B${row2}").Value2
    $rest1 = $ws.Range("E${row1}:AD${row1}").Value2
    $rest2 = $ws.Range("E${row2}:AD${row2}").Value2

    $ws.Range("B${row1}").Value = $b2
    $ws.Range("B${row2}").Value = $b1

    $ws.Range("E${row1}:AD${row1}").Value = $rest2
    $ws.Range("E${row2}:AD${row2}").Value = $rest1
}

Swap-Rows 110 111
Swap-Rows 237 238
